# Updated cryptos list with GitHub Actions
# Applies the latest price / 1h-volume-change refresh pulled from
# coinranking.com, including the OKB <-> FLOKI rank swap (rows 50-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.538.22'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '3.096.92'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '234.92'
$ws.Range("E5").Value = '  +9.59%  '
$ws.Range("D6").Value = '625.68'
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '3.095.84'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").Value = '0.726'
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  +5.81%  '
$ws.Range("D14").Value = '36.65'
$ws.Range("E14").Value = '  +6.16%  '
$ws.Range("D15").Value = '5.48'
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").Value = '90.078.92'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '3.082.77'
$ws.Range("E18").Value = '  -3.23%  '
$ws.Range("D19").Value = '3.82'
$ws.Range("E19").Value = '  +5.29%  '
$ws.Range("D20").Value = '0.0000218'
$ws.Range("E20").Value = '  +4.87%  '
$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").Value = '437.88'
$ws.Range("E22").Value = '  -3.30%  '
$ws.Range("D23").Value = '5.55'
$ws.Range("E23").Value = '  +6.32%  '
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").Value = '5.94'
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("D26").Value = '7.58'
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("D27").Value = '88.63'
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("D28").Value = '12.19'
$ws.Range("E28").Value = '  +2.59%  '
$ws.Range("D29").Value = '3.253.59'
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '9.44'
$ws.Range("E31").Value = '  +4.39%  '
# D32 is a trailing-zero decimal ("0.160"); force Text format so Excel
# doesn't normalize it down to "0.16" the way it would a plain numeric Value.
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.160'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '0.195'
$ws.Range("E33").Value = '  +11.83%  '
$ws.Range("D34").Value = '3.89'
$ws.Range("E34").Value = '  +8.21%  '
$ws.Range("E35").Value = '  +8.79%  '
# D36 likewise ("25.80" has a significant trailing zero).
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '25.80'
$ws.Range("E36").Value = '  -2.55%  '
$ws.Range("D37").Value = '0.893'
$ws.Range("E37").Value = '  -10.96%  '
$ws.Range("D38").Value = '507.78'
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("E39").Value = '  +5.54%  '
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("E41").Value = '  +2.86%  '
$ws.Range("E42").Value = '  -0.08%  '
# D43 likewise ("0.410" has a significant trailing zero).
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.410'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '3.46'
$ws.Range("E46").Value = '  +55.85%  '
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("D48").Value = '151.86'
$ws.Range("E48").Value = '  +3.49%  '
$ws.Range("D49").Value = '0.688'
$ws.Range("E49").Value = '  +5.93%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '0.000278'
$ws.Range("E50").Value = '  +19.06%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '44.93'
$ws.Range("E51").Value = '  +1.31%  '
